# Insert a new data row before the existing row 928 (shifting rows 928:986
# down to 929:987) and populate it with the new "Patagonia" Papa entry for
# Vega Modelo de Temuco (La Araucanía).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(928).Insert()

$ws.Range("A928").Value = 10
$ws.Range("B928").Value = "Vega Modelo de Temuco"
$ws.Range("C928").Value = "La Araucanía"
$ws.Range("D928").Value = 44931
$ws.Range("E928").Value = 9
$ws.Range("F928").Value = 100114001
$ws.Range("G928").Value = "Papa"
$ws.Range("H928").Value = "Patagonia"
$ws.Range("I928").Value = "1a nueva(o)"
$ws.Range("J928").Value = 1500
$ws.Range("K928").Value = 11000
$ws.Range("L928").Value = 12000
$ws.Range("M928").Value = 11467
$ws.Range("N928").Value = "`$/saco 25 kilos"
$ws.Range("O928").Value = "Provincia de Cautín"
$ws.Range("P928").Value = 459
$ws.Range("Q928").Value = 25
$ws.Range("R928").Value = "Hortaliza"
